$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Runtime without error -> "no", with note explaining why
$ws.Range("B6").Value = "no"
$ws.Range("C6").Value = "It doesn't redirect to correct page"

# Assertion validity -> clear value and note (now blank)
$ws.Range("B7").Value = $null
$ws.Range("C7").Value = $null

# Code BLEU updated score + note text
$ws.Range("B12").Value = 0.3065697598209665
$ws.Range("C12").Value = "{'codebleu': 0.30656975982096646, 'ngram_match_score': 0.17709940898665436, 'weighted_ngram_match_score': 0.22650915009248465, 'syntax_match_score': 0.547945205479452, 'dataflow_match_score': 0.27472527472527475}"

# Update the active selection to match the saved view state
$ws.Range("B7").Select()
